$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - overwrite the existing row (was ID/username/password/email) with
# the first user's data, left to right.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = "password"
$ws.Range("D2").Value = "Shakira Regalado"
$ws.Range("E2").Value = "shakiraregalado@gmail.com"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "Fluffy"

# Row 3 - second user's data, left to right.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "kira"
$ws.Range("C3").Value = "pass"
$ws.Range("D3").Value = "Shakira"
$ws.Range("E3").Value = "shakira@gmail.com"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = "Pink"

# Row 1 - headers added last, left to right.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Full Name"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Security Question"
$ws.Range("G1").Value = "Security Answer"

# Column widths as set by Excel after data entry / autofit
$ws.Columns.Item(5).ColumnWidth = 18.3
$ws.Columns.Item(6).ColumnWidth = 18

$ws.Range("G2").Select()
